$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Snapshot the existing "highlight" (yellow) cell format on C2 and paste
# it onto E3 *before* we touch C2's own formatting, so the yellow fill
# moves from C2 to E3 without creating any duplicate style/fill entries. ---
$ws.Range("C2").Copy()
$ws.Range("E3").PasteSpecial(-4122)   # xlPasteFormats

# Re-base C2's formatting on a plain (non-highlighted) cell that already
# uses the "centered / wrapped, no fill" style (e.g. B1), so C2 loses its
# yellow highlight but keeps the centered/wrapped look.
$ws.Range("B1").Copy()
$ws.Range("C2").PasteSpecial(-4122)   # xlPasteFormats

$excel.CutCopyMode = $false

# --- Cell value updates -----------------------------------------------
$ws.Range("C2").Value = "سفارشات"
$ws.Range("A3").Value = "لاگین از گوگل"
$ws.Range("C3").Value = "سبد خرید عملیاتی"
$ws.Range("E3").Value = "تصاویر"
$ws.Range("E4").Value = "لاگین"
$ws.Range("E5").Value = "تغییر رمز پروفایل"
$ws.Range("E6").Value = "پنل کاربری"
$ws.Range("E7").Value = "آدرس ها"
$ws.Range("E8").Value = "سبد خرید"
$ws.Range("E9").Value = "ثبت نام"
$ws.Range("E10").Value = "پنل ادمین"
$ws.Range("B13").Value = "confirm email and phone number"
$ws.Range("B14").Value = "ارسال اس ام اس"
$ws.Range("B15").Value = "ارسال ایمیل"

# --- Cells that no longer hold content (drop the cell entirely, like the
# other always-empty cells in this sheet, rather than leaving a blank
# styled cell behind) ---------------------------------------------------
$ws.Range("A4").Clear()
$ws.Range("C4").Clear()
$ws.Range("C5").Clear()

# --- Row 13 grew taller to fit the new wrapped text ---------------------
$ws.Rows(13).RowHeight = 60

# --- Selection moved ------------------------------------------------
[void]$ws.Range("J6").Select()
